$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at the very top (title row) -- shifts old row1 (headers) to row2, etc.
$ws.Rows.Item(1).Insert()

# 2. Insert a new row after the header row (now row2) to hold the new BHVO 2G standard sample.
$ws.Rows.Item(3).Insert()
$ws.Range("M3").Clear()

# 3. New BHVO 2G standard row (row 3) -- no Total H2O value for the standard.
$ws.Range("A3").Value = "BHVO 2G"
$ws.Range("B3").Value = "Standard"
$ws.Range("C3").Value = 49.3
$ws.Range("D3").Value = 2.79
$ws.Range("E3").Value = 13.6
$ws.Range("F3").Value = 11.3
$ws.Range("G3").Value = 0.17
$ws.Range("H3").Value = 7.13
$ws.Range("I3").Value = 11.4
$ws.Range("J3").Value = 2.4
$ws.Range("K3").Value = 0.51
$ws.Range("L3").Value = 98.6

# 4. Footnote row at the bottom (row 48, leaving row 47 blank)
$ws.Range("A48").Value = "BHVO 2G: Standard sample (Coulthard Jr 2018); Sid: Sideromelane; Pal: Palagonite"

# 5. Title in A1 (added last so it lands at the highest shared-string index)
$ws.Range("A1").Value = "Table S1 Major elements of tested sideromelanes and palagonites"

Write-Host "done"
